# Reformat excel template and add up to 5 answers
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Survey")

# Start from a clean slate for the data region (drops old values + formats)
$ws.Range("A1:H4").Clear()

# ---- Header row ----
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Question"
$ws.Range("C1").Value = "Solution"
$ws.Range("D1").Value = "StartSet"
$ws.Range("E1").Value = "Difficulty"
$ws.Range("F1").Value = "Slope"
$ws.Range("G1").Value = "A1"
$ws.Range("H1").Value = "A2"
$ws.Range("I1").Value = "A3"
$ws.Range("J1").Value = "A4"
$ws.Range("K1").Value = "A5"

# ---- Row 2 ----
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Pizza?"
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = "X"
$ws.Range("E2").Value = 1.2
$ws.Range("F2").Value = 10.1
$ws.Range("F2").NumberFormat = "General"
$ws.Range("G2").Value = "Tomaten"
$ws.Range("H2").Value = "Schinken"
$ws.Range("I2").Value = 3
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 5

# ---- Row 3 ----
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Döner?"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = "X"
$ws.Range("E3").Value = 0.5
$ws.Range("F3").Value = 1
$ws.Range("F3").NumberFormat = "0.00"
$ws.Range("G3").Value = "Scharf"
$ws.Range("H3").Value = "Käse"
$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 4
$ws.Range("K3").Value = 5

# ---- Row 4 ----
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Nudeln?"
$ws.Range("C4").Value = 2
$ws.Range("E4").Value = 1.1000000000000001
$ws.Range("F4").Value = 0.75
$ws.Range("F4").NumberFormat = "General"
$ws.Range("G4").Value = "Tomaten"
$ws.Range("H4").Value = "Sahne"
$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 4
$ws.Range("K4").Value = 5

# ---- Column widths ----
$ws.Range("A1").EntireColumn.ColumnWidth = 2.85546875
$ws.Range("B1:C1").EntireColumn.ColumnWidth = 12.7109375
$ws.Range("D1").EntireColumn.ColumnWidth = 13
$ws.Range("E1:H1").EntireColumn.ColumnWidth = 12.7109375

# ---- Selection ----
$ws.Range("F8").Select()
